$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a leading apostrophe to force Excel to store these as literal text
# strings (matching the original inlineStr/text cells) instead of letting
# Excel auto-convert numeric-looking text into numbers or percentages.
$ws.Range("D2").Value = "'286.07"
$ws.Range("E2").Value = "'3.67%"
$ws.Range("D3").Value = "'28.75"
$ws.Range("E3").Value = "'5.92%"
$ws.Range("D4").Value = "'4.928"
$ws.Range("E4").Value = "'1.37%"
$ws.Range("D5").Value = "'0.06528"
$ws.Range("E5").Value = "'1.98%"
$ws.Range("D6").Value = "'7.241"
$ws.Range("E6").Value = "'4.23%"
$ws.Range("D7").Value = "'1.363"
$ws.Range("E7").Value = "'14.24%"
$ws.Range("D8").Value = "'0.9110"
$ws.Range("E8").Value = "'3.92%"
$ws.Range("D9").Value = "'0.1555"
$ws.Range("E9").Value = "'2.59%"
$ws.Range("D10").Value = "'0.06819"
$ws.Range("E10").Value = "'34.07%"
$ws.Range("D11").Value = "'0.07693"
$ws.Range("E11").Value = "'2.18%"
$ws.Range("D12").Value = "'0.02995"
$ws.Range("E12").Value = "'0.91%"
$ws.Range("D13").Value = "'0.08969"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("D14").Value = "'0.001599"
$ws.Range("E14").Value = "'2.40%"
$ws.Range("D15").Value = "'0.0006531"
$ws.Range("E15").Value = "'2.29%"
$ws.Range("D16").Value = "'0.006029"
$ws.Range("E16").Value = "'-2.42%"
$ws.Range("D17").Value = "'3.462"
$ws.Range("E17").Value = "'-0.58%"
$ws.Range("D18").Value = "'3.389"
$ws.Range("E18").Value = "'2.46%"
$ws.Range("D19").Value = "'2.242"
$ws.Range("E19").Value = "'-0.46%"
$ws.Range("D20").Value = "'0.3150"
$ws.Range("E20").Value = "'0.44%"
$ws.Range("D21").Value = "'0.1350"
$ws.Range("E21").Value = "'0.08%"
$ws.Range("D22").Value = "'4.048"
$ws.Range("E22").Value = "'3.33%"
$ws.Range("D23").Value = "'0.1555"
$ws.Range("E23").Value = "'12.66%"
$ws.Range("D24").Value = "'0.04470"
$ws.Range("E24").Value = "'1.35%"
$ws.Range("D25").Value = "'0.001190"
$ws.Range("E25").Value = "'1.15%"
$ws.Range("D26").Value = "'0.004328"
$ws.Range("E26").Value = "'12.03%"
$ws.Range("D28").Value = "'0.0001182"
$ws.Range("E28").Value = "'-1.47%"
$ws.Range("E29").Value = "'-15.61%"
$ws.Range("E40").Value = "'0.26%"
$ws.Range("D41").Value = "'0.006695"
$ws.Range("E41").Value = "'-1.28%"
$ws.Range("E42").Value = "'5.21%"
$ws.Range("D43").Value = "'0.002163"
$ws.Range("E43").Value = "'0.20%"
$ws.Range("D44").Value = "'0.01176"
$ws.Range("E44").Value = "'-0.52%"
$ws.Range("D45").Value = "'0.00005407"
$ws.Range("E45").Value = "'2.83%"
$ws.Range("E46").Value = "'-7.47%"
$ws.Range("D47").Value = "'0.01853"
$ws.Range("E47").Value = "'0.20%"
